{"js": "// Replace each \"AxB=\" multiplication prompt in the table with its new\n// value, per the commit's regenerated number set. Every old prompt string\n// is unique in the document, so a direct search + whole-text replace on\n// each matching run is safe (no cross-matches between old/new sets).\nconst replacements = [\n  [\"70\u00d797=\", \"45\u00d786=\"],\n  [\"38\u00d761=\", \"89\u00d723=\"],\n  [\"73\u00d736=\", \"98\u00d740=\"],\n  [\"74\u00d770=\", \"54\u00d740=\"],\n  [\"58\u00d751=\", \"20\u00d752=\"],\n  [\"11\u00d794=\", \"23\u00d713=\"],\n  [\"97\u00d727=\", \"71\u00d785=\"],\n  [\"38\u00d746=\", \"28\u00d721=\"],\n  [\"23\u00d785=\", \"78\u00d750=\"],\n  [\"62\u00d713=\", \"47\u00d777=\"],\n  [\"47\u00d780=\", \"47\u00d734=\"],\n  [\"57\u00d781=\", \"43\u00d732=\"],\n  [\"52\u00d790=\", \"94\u00d739=\"],\n  [\"47\u00d787=\", \"53\u00d798=\"],\n  [\"38\u00d714=\", \"38\u00d760=\"],\n  [\"74\u00d758=\", \"85\u00d763=\"],\n  [\"29\u00d797=\", \"34\u00d732=\"],\n  [\"16\u00d728=\", \"65\u00d713=\"],\n  [\"29\u00d776=\", \"30\u00d744=\"],\n  [\"67\u00d736=\", \"33\u00d795=\"],\n  [\"18\u00d723=\", \"67\u00d749=\"],\n  [\"98\u00d743=\", \"64\u00d725=\"],\n  [\"61\u00d738=\", \"76\u00d780=\"],\n  [\"95\u00d790=\", \"47\u00d732=\"],\n  [\"28\u00d775=\", \"97\u00d716=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items,text\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"AxB=\" multiplication prompt in the table with its new\n# value, per the commit's regenerated number set. Every old prompt string\n# is unique in the document, so Find/Replace per pair is safe (no\n# cross-matches between the old and new sets).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"70\u00d797=\", \"45\u00d786=\"),\n    @(\"38\u00d761=\", \"89\u00d723=\"),\n    @(\"73\u00d736=\", \"98\u00d740=\"),\n    @(\"74\u00d770=\", \"54\u00d740=\"),\n    @(\"58\u00d751=\", \"20\u00d752=\"),\n    @(\"11\u00d794=\", \"23\u00d713=\"),\n    @(\"97\u00d727=\", \"71\u00d785=\"),\n    @(\"38\u00d746=\", \"28\u00d721=\"),\n    @(\"23\u00d785=\", \"78\u00d750=\"),\n    @(\"62\u00d713=\", \"47\u00d777=\"),\n    @(\"47\u00d780=\", \"47\u00d734=\"),\n    @(\"57\u00d781=\", \"43\u00d732=\"),\n    @(\"52\u00d790=\", \"94\u00d739=\"),\n    @(\"47\u00d787=\", \"53\u00d798=\"),\n    @(\"38\u00d714=\", \"38\u00d760=\"),\n    @(\"74\u00d758=\", \"85\u00d763=\"),\n    @(\"29\u00d797=\", \"34\u00d732=\"),\n    @(\"16\u00d728=\", \"65\u00d713=\"),\n    @(\"29\u00d776=\", \"30\u00d744=\"),\n    @(\"67\u00d736=\", \"33\u00d795=\"),\n    @(\"18\u00d723=\", \"67\u00d749=\"),\n    @(\"98\u00d743=\", \"64\u00d725=\"),\n    @(\"61\u00d738=\", \"76\u00d780=\"),\n    @(\"95\u00d790=\", \"47\u00d732=\"),\n    @(\"28\u00d775=\", \"97\u00d716=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $oldText,    # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $null, $null, $true, 1, $true,\n        $newText,    # ReplaceWith\n        2            # wdReplaceAll\n    )\n}\n"}
